$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "admin123"
$ws.Range("B2").Value = "Admin"
